# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Widen the "Error Detail" column (P) on both locale sheets to 40 characters.
# (39.1666... is used as the COM ColumnWidth input because the host's
# character-width<->pixel round trip adds ~5/6 of a character back on
# save; this lands exactly on a saved width of 40, matching the other
# already-40-wide columns in this workbook.)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664

# Record the handback-transform failure + detail message for each locale.
$wsZhCn.Range("P3").Value = "Handback file name: wbeuve4z.i3i is different with handoff file name: 30f92c11-1db5-467e-93b9-725b4fa15fde.f341e554cef27a70f1feea1ad5ad2af7e6f9786b.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: wbeuve4z.i3i is different with handoff file name: 30f92c11-1db5-467e-93b9-725b4fa15fde.f341e554cef27a70f1feea1ad5ad2af7e6f9786b.de-de."

# Update the shared "Status" text used across Overview / zh-cn / de-de sheets.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"
